$d = $word.ActiveDocument

$replacements = @(
    @("2025-06-08 Sunday", "2025-06-09 Monday"),
    @("91÷9=", "96÷2="),
    @("68÷6=", "14÷2="),
    @("81÷5=", "36÷3="),
    @("82÷4=", "84÷4="),
    @("51÷9=", "31÷2="),
    @("53÷5=", "63÷2="),
    @("10÷5=", "86÷5="),
    @("92÷6=", "79÷6="),
    @("15÷2=", "31÷7="),
    @("92÷2=", "19÷7="),
    @("96÷5=", "39÷2="),
    @("22÷4=", "67÷3="),
    @("37÷4=", "40÷8="),
    @("12÷2=", "88÷6="),
    @("69÷3=", "48÷4="),
    @("17÷8=", "96÷3="),
    @("90÷6=", "73÷9="),
    @("33÷7=", "87÷6="),
    @("56÷5=", "37÷7="),
    @("23÷2=", "30÷6="),
    @("22÷2=", "59÷3="),
    @("47÷2=", "11÷5="),
    @("97÷2=", "81÷3="),
    @("26÷4=", "76÷4="),
    @("46÷5=", "34÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
